$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# --- Shrink the "Status" columns to match the shorter text ---
# (ColumnWidth is expressed in characters; Excel snaps the stored
# worksheet width to its internal pixel grid, so 12.5 is the value
# that lands closest to the target column width.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
